# Anonymize feature analysis counts: zero out the per-user email counts
# on the SentEmailCounts and ReceivedEmailCounts sheets (column B, rows 2-147).
# Any charts referencing these ranges will refresh their cached values
# automatically when the workbook is saved/recalculated.

$wb = $excel.ActiveWorkbook

$sheetNames = @("SentEmailCounts", "ReceivedEmailCounts")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $rng = $ws.Range("B2:B147")
    $zeros = New-Object 'object[,]' 146,1
    for ($i = 0; $i -lt 146; $i++) {
        $zeros[$i, 0] = 0
    }
    $rng.Value = $zeros
}

$excel.CalculateFullRebuild()
